# "Generate Report for Handback"
#
# The localization-status report is refreshed with handback information:
#  - the "Ready for handoff" status text becomes "Handed back: in sync
#    with en-US" everywhere it appears (Overview sheet Status columns and
#    each per-locale sheet's Status column),
#  - two new columns of data are populated on the per-locale sheets:
#    "Latest Target File" (F) and "Latest Handback File" (G) — these
#    mirror the existing handoff file name (markdown source + the
#    per-locale xlf) with their own hyperlinks,
#  - the "Latest Handback DateTime" (H) column, previously the zero-date
#    placeholder, gets a real handback timestamp (one shared timestamp
#    per locale sheet).

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Update the "Status" text wherever it is used ---------------------

$wsOverview.Range("B2").Value = $statusNew
$wsOverview.Range("C2").Value = $statusNew
$wsOverview.Range("B3").Value = $statusNew
$wsOverview.Range("C3").Value = $statusNew

$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("C3").Value = $statusNew

$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("C3").Value = $statusNew

# --- 2. zh-cn sheet: populate Latest Target File / Latest Handback File --

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3a68918dc30613db901c44b80c8cd3546a376665/e2e/cac2d91b-a7d4-4b7c-bcf8-598dce9987e5.md",
    "",
    "",
    "cac2d91b-a7d4-4b7c-bcf8-598dce9987e5.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c6bcc9b49ba688c7d5fd93f30572bb51729b0c9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cac2d91b-a7d4-4b7c-bcf8-598dce9987e5.27bb06dd54e7b2ba5573bb897946066b1624666c.zh-cn.xlf",
    "",
    "",
    "cac2d91b-a7d4-4b7c-bcf8-598dce9987e5.27bb06dd54e7b2ba5573bb897946066b1624666c.zh-cn.xlf"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3a68918dc30613db901c44b80c8cd3546a376665/e2e/f7c14317-9349-4704-94f9-d74396d8d4cf.md",
    "",
    "",
    "f7c14317-9349-4704-94f9-d74396d8d4cf.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c6bcc9b49ba688c7d5fd93f30572bb51729b0c9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f7c14317-9349-4704-94f9-d74396d8d4cf.7f346e2eddbe90e5f91da9d35de42219eed3788f.zh-cn.xlf",
    "",
    "",
    "f7c14317-9349-4704-94f9-d74396d8d4cf.7f346e2eddbe90e5f91da9d35de42219eed3788f.zh-cn.xlf"
)

# zh-cn handback timestamp (shared by both data rows)
$wsZhCn.Range("H2").Value = "2016-03-22 18:49:04"
$wsZhCn.Range("H3").Value = "2016-03-22 18:49:04"

# --- 3. de-de sheet: populate Latest Target File / Latest Handback File --

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3a68918dc30613db901c44b80c8cd3546a376665/e2e/cac2d91b-a7d4-4b7c-bcf8-598dce9987e5.md",
    "",
    "",
    "cac2d91b-a7d4-4b7c-bcf8-598dce9987e5.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea6801aecb92a4969c526d54f7a14b6c910fa640/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cac2d91b-a7d4-4b7c-bcf8-598dce9987e5.27bb06dd54e7b2ba5573bb897946066b1624666c.de-de.xlf",
    "",
    "",
    "cac2d91b-a7d4-4b7c-bcf8-598dce9987e5.27bb06dd54e7b2ba5573bb897946066b1624666c.de-de.xlf"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3a68918dc30613db901c44b80c8cd3546a376665/e2e/f7c14317-9349-4704-94f9-d74396d8d4cf.md",
    "",
    "",
    "f7c14317-9349-4704-94f9-d74396d8d4cf.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea6801aecb92a4969c526d54f7a14b6c910fa640/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f7c14317-9349-4704-94f9-d74396d8d4cf.7f346e2eddbe90e5f91da9d35de42219eed3788f.de-de.xlf",
    "",
    "",
    "f7c14317-9349-4704-94f9-d74396d8d4cf.7f346e2eddbe90e5f91da9d35de42219eed3788f.de-de.xlf"
)

# de-de handback timestamp (shared by both data rows, distinct from zh-cn's)
$wsDeDe.Range("H2").Value = "2016-03-22 18:49:12"
$wsDeDe.Range("H3").Value = "2016-03-22 18:49:12"
